function Set-CellText($sheet, $ref, $text) {
    $cell = $sheet.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-CellText $ws 'D2' '61.861.10'
Set-CellText $ws 'E2' '  +3.62%  '
Set-CellText $ws 'D3' '3.407.38'
Set-CellText $ws 'E3' '  +2.30%  '
Set-CellText $ws 'E4' '  +0.09%  '
Set-CellText $ws 'D5' '577.34'
Set-CellText $ws 'E5' '  +2.67%  '
Set-CellText $ws 'D6' '137.61'
Set-CellText $ws 'E6' '  +5.93%  '
Set-CellText $ws 'E7' '  +0.06%  '
Set-CellText $ws 'D8' '3.406.56'
Set-CellText $ws 'E8' '  +2.30%  '
Set-CellText $ws 'D9' '0.477'
Set-CellText $ws 'E9' '  +1.10%  '
Set-CellText $ws 'D10' '7.49'
Set-CellText $ws 'E10' '  +0.91%  '
Set-CellText $ws 'E11' '  +7.87%  '
Set-CellText $ws 'E12' '  +5.11%  '
Set-CellText $ws 'D13' '3.994.73'
Set-CellText $ws 'E13' '  +2.78%  '
Set-CellText $ws 'E14' '  +2.18%  '
Set-CellText $ws 'E15' '  +6.43%  '
Set-CellText $ws 'D16' '3.409.79'
Set-CellText $ws 'E16' '  +2.92%  '
Set-CellText $ws 'D17' '25.42'
Set-CellText $ws 'E17' '  +4.42%  '
Set-CellText $ws 'D18' '61.858.54'
Set-CellText $ws 'E18' '  +3.31%  '
Set-CellText $ws 'D19' '14.23'
Set-CellText $ws 'E19' '  +6.59%  '
Set-CellText $ws 'E20' '  +4.48%  '
Set-CellText $ws 'D21' '9.51'
Set-CellText $ws 'E21' '  +5.33%  '
Set-CellText $ws 'D22' '389.73'
Set-CellText $ws 'E22' '  +10.52%  '
Set-CellText $ws 'D23' '0.572'
Set-CellText $ws 'E23' '  +3.07%  '
Set-CellText $ws 'D24' '3.544.83'
Set-CellText $ws 'E24' '  +2.61%  '
Set-CellText $ws 'E25' '  +15.71%  '
Set-CellText $ws 'E26' '  +0.00%  '
Set-CellText $ws 'D27' '71.54'
Set-CellText $ws 'E27' '  +3.91%  '
Set-CellText $ws 'E28' '  +3.78%  '
Set-CellText $ws 'D29' '1.59'
Set-CellText $ws 'E29' '  +5.36%  '
Set-CellText $ws 'D30' '0.998'
Set-CellText $ws 'E30' '  +0.00%  '
Set-CellText $ws 'E31' '  +5.31%  '
Set-CellText $ws 'E32' '  +5.35%  '
Set-CellText $ws 'E33' '  +3.17%  '
Set-CellText $ws 'D34' '3.438.27'
Set-CellText $ws 'E34' '  +2.42%  '
Set-CellText $ws 'E35' '  +0.02%  '
Set-CellText $ws 'D36' '23.55'
Set-CellText $ws 'E36' '  +3.30%  '
Set-CellText $ws 'D37' '5.44'
Set-CellText $ws 'E37' '  +2.47%  '
Set-CellText $ws 'E38' '  +2.79%  '
Set-CellText $ws 'E39' '  +4.52%  '
Set-CellText $ws 'D40' '164.61'
Set-CellText $ws 'E40' '  +4.16%  '
Set-CellText $ws 'D41' '0.0789'
Set-CellText $ws 'E41' '  +4.62%  '
Set-CellText $ws 'E42' '  +13.78%  '
Set-CellText $ws 'D43' '0.788'
Set-CellText $ws 'E43' '  +5.70%  '
Set-CellText $ws 'E44' '  +0.05%  '
Set-CellText $ws 'E45' '  +3.23%  '
Set-CellText $ws 'E46' '  +2.34%  '
Set-CellText $ws 'D47' '41.66'
Set-CellText $ws 'E47' '  +2.02%  '
Set-CellText $ws 'D48' '24.83'
Set-CellText $ws 'E48' '  +7.12%  '
Set-CellText $ws 'D49' '6.98'
Set-CellText $ws 'E49' '  +2.89%  '
Set-CellText $ws 'D50' '23.25'
Set-CellText $ws 'E50' '  +4.44%  '
Set-CellText $ws 'D51' '2.381.70'
Set-CellText $ws 'E51' '  +10.68%  '
